$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.106.10"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "1.638.17"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.47"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.520"
$ws.Range("E6").Value = "  +1.56%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.97"
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").Value = "1.640.07"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.540"
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.63"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("D17").Value = "27.135.01"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.81"
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.95"
$ws.Range("E21").Value = "  +1.84%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.41"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.51"
$ws.Range("E23").Value = "  +3.24%  "
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.71"
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.41"
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.67"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.37"
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("D34").Value = "1.306.08"
$ws.Range("E34").Value = "  +2.90%  "
$ws.Range("E36").Value = "  +1.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0175"
$ws.Range("E37").Value = "  -1.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.856"
$ws.Range("E38").Value = "  +2.90%  "
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.812"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("E42").Value = "  +5.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.29"
$ws.Range("E43").Value = "  -1.81%  "
$ws.Range("D44").Value = "1.778.13"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.80"
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.35"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("E48").Value = "  +2.28%  "
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.62"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0961"
$ws.Range("E51").Value = "  -0.02%  "
